$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the visit date values in I4:I7 from 46007 (2025-12-16) to 46009 (2025-12-18)
$ws.Range("I4:I7").Value = 46009

# Update the active cell / selection to I16
$ws.Range("I16").Select()
